# Updated cryptos list values (price & 1h volume change) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.010.30"
$ws.Range("E2").Value = "'  +1.41%  "

$ws.Range("D3").Value = "'3.307.91"
$ws.Range("E3").Value = "'  +5.54%  "

$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'598.74"
$ws.Range("E5").Value = "'  +0.33%  "

$ws.Range("D6").Value = "'143.18"
$ws.Range("E6").Value = "'  +4.48%  "

$ws.Range("E7").Value = "'  +0.07%  "

$ws.Range("D8").Value = "'3.303.75"
$ws.Range("E8").Value = "'  +5.84%  "

$ws.Range("E9").Value = "'  +0.46%  "

$ws.Range("E10").Value = "'  +2.14%  "

$ws.Range("D11").Value = "'5.46"
$ws.Range("E11").Value = "'  +3.01%  "

$ws.Range("D12").Value = "'0.472"
$ws.Range("E12").Value = "'  +2.36%  "

$ws.Range("E13").Value = "'  -0.39%  "

$ws.Range("D14").Value = "'34.88"
$ws.Range("E14").Value = "'  +1.50%  "

$ws.Range("D15").Value = "'3.850.76"
$ws.Range("E15").Value = "'  +5.80%  "

$ws.Range("E16").Value = "'  +1.19%  "

$ws.Range("D17").Value = "'3.308.28"

$ws.Range("D18").Value = "'64.114.44"
$ws.Range("E18").Value = "'  +1.46%  "

$ws.Range("D19").Value = "'6.89"
$ws.Range("E19").Value = "'  +2.20%  "

$ws.Range("D20").Value = "'483.28"
$ws.Range("E20").Value = "'  +1.14%  "

$ws.Range("D21").Value = "'14.29"
$ws.Range("E21").Value = "'  +0.77%  "

$ws.Range("E22").Value = "'  +6.20%  "

$ws.Range("D23").Value = "'8.02"
$ws.Range("E23").Value = "'  +3.72%  "

$ws.Range("D24").Value = "'13.52"
$ws.Range("E24").Value = "'  +3.38%  "

$ws.Range("D25").Value = "'84.41"
$ws.Range("E25").Value = "'  -3.33%  "

$ws.Range("E26").Value = "'  -0.04%  "

$ws.Range("E27").Value = "'  +2.09%  "

$ws.Range("E28").Value = "'  +2.31%  "

$ws.Range("D29").Value = "'8.28"
$ws.Range("E29").Value = "'  +3.24%  "

$ws.Range("E30").Value = "'  +0.00%  "

$ws.Range("E31").Value = "'  +2.22%  "

$ws.Range("D32").Value = "'28.68"
$ws.Range("E32").Value = "'  +5.70%  "

$ws.Range("E33").Value = "'  -1.45%  "

$ws.Range("E34").Value = "'  +0.66%  "

$ws.Range("E35").Value = "'  +1.53%  "

$ws.Range("D36").Value = "'5.99"
$ws.Range("E36").Value = "'  +2.30%  "

$ws.Range("D37").Value = "'53.41"
$ws.Range("E37").Value = "'  +2.75%  "

$ws.Range("D38").Value = "'0.0₃0738"

$ws.Range("E39").Value = "'  +2.36%  "

$ws.Range("D40").Value = "'432.45"
$ws.Range("E40").Value = "'  +2.82%  "

$ws.Range("D41").Value = "'3.016.37"
$ws.Range("E41").Value = "'  +4.50%  "

$ws.Range("D42").Value = "'8.45"
$ws.Range("E42").Value = "'  +1.87%  "

$ws.Range("D43").Value = "'2.77"
$ws.Range("E43").Value = "'  +2.91%  "

$ws.Range("E44").Value = "'  -5.80%  "

$ws.Range("E45").Value = "'  +1.99%  "

$ws.Range("E46").Value = "'  +4.58%  "

$ws.Range("D47").Value = "'26.30"
$ws.Range("E47").Value = "'  +2.13%  "

$ws.Range("E48").Value = "'  +0.07%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'2.33"
$ws.Range("E49").Value = "'  +2.10%  "

$ws.Range("E50").Value = "'  +1.51%  "

$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'35.43"
$ws.Range("E51").Value = "'  +14.38%  "
